# Apply the requested edits to the Error_types_table worksheet:
#  - Column B (Query): replace long filenames with short human-readable labels
#  - Column D (Outcome): replace "Genus + species ..." with "Genus and species ..."

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of old filename substrings -> new short labels
$queryMap = @{
    "make_12s_16s_simulated_reads_5-BetterDatabaseARTSimulation_runEDNAFLOW_12S_Lulu_RESULTS_dada2_asv.fa" = "100 Australian species"
    "make_12s_16s_simulated_reads_5-BetterDatabaseARTSimulation_runEDNAFLOW_16S_Lulu_RESULTS_dada2_asv.fa" = "100 Australian species"
    "make_12s_16s_simulated_reads_7-Lutjanids_Mock_runEDNAFlow_12S_RESULTS_dada2_asv.fa" = "Lutjanidae"
    "make_12s_16s_simulated_reads_7-Lutjanids_Mock_runEDNAFlow_16S_RESULTS_dada2_asv.fa" = "Lutjanidae"
    "make_12s_16s_simulated_reads_8-Rottnest_runEDNAFLOW_12S_RESULTS_dada2_asv.fa" = "Rottnest"
    "make_12s_16s_simulated_reads_8-Rottnest_runEDNAFLOW_16S_RESULTS_dada2_asv.fa" = "Rottnest"
}

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $bCell = $ws.Cells.Item($r, 2)
    $bVal = $bCell.Value2
    if ($null -ne $bVal -and $queryMap.ContainsKey($bVal)) {
        $bCell.Value2 = $queryMap[$bVal]
    }

    $dCell = $ws.Cells.Item($r, 4)
    $dVal = $dCell.Value2
    if ($dVal -eq "Genus + species correct") {
        $dCell.Value2 = "Genus and species correct"
    } elseif ($dVal -eq "Genus + species wrong") {
        $dCell.Value2 = "Genus and species wrong"
    }
}
